# Features List.xlsx update
# - Updates the Garrett ("Dev") task list (rows 4-15) with newly completed /
#   added features and removes the now-obsolete "block quality" & "closest
#   player" placeholder rows.
# - Updates the Tobi task list (rows 19-23): reworded bug/feature rows and
#   blanks out the row that is no longer populated.
# - Widens column B slightly and moves the active selection to the new last
#   row (B15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($rowNum, $b, $c, $d, $wrapD) {
    $ws.Cells.Item($rowNum, 2).Value = $b
    if ($null -ne $c) { $ws.Cells.Item($rowNum, 3).Value = $c }
    if ($null -ne $d) {
        $ws.Cells.Item($rowNum, 4).Value = $d
        if ($wrapD) { $ws.Cells.Item($rowNum, 4).WrapText = $true }
    }
}

# ---- Garrett section (rows 4-15), strikethrough "done" style already set ----
Set-Row 4  "When the set choice is made, ball is sent to the players location" "Garrett" "When a set choice is made by either player or AI, send the ball to the location of that player, and not just a tile on the grid"
Set-Row 5  "AI serve location determined" "Garrett" "The ai randomly determines a serve location on the court"
Set-Row 6  "AI serve location displayed by the ball's position" "Garrett" "The AI's serve location is displayed by the ball position on the court map"
Set-Row 7  "Pass/Dig ball position based on quality" "Garrett" "The ball position on the map is placed in one of three locations based on the quality of the pass (1, 2 and 3 passes) for both player and AI"
Set-Row 8  "AI will move into defensive positions" "Garrett" "At each stage of the rally, the ai will move into appropriate positions on the court"
Set-Row 9  "Ball serve location based on serve quality" "Garrett" "When the player or AI serves the ball, if the quality is low the ball position can be changed randomly"
Set-Row 10 "AI chooses where to attack" "Garrett" "Have the AI choose randomly where to attack, to be improved in future"
Set-Row 11 "Ball attack location change based on quality" "Garrett" "The location of the ball can change randomly based on the quality of the attack, if the attack randomly goes off the grid, it is an error. Need to also reduce the chances for hitting errors to account for these new error chances"

# Row 12 is new - carry forward the same "done" strikethrough style used by rows 4-11
$ws.Range("B12:D12").Font.Strikethrough = $true
Set-Row 12 "Function for finding closest player to the ball on attacks and defence" "Garrett" "Done"

# Row 13 - new row, only a feature name (no dev/details yet)
$ws.Range("C13:D13").Font.Strikethrough = $false
$ws.Range("C13:D13").Value = $null
Set-Row 13 "Edge Case: setter makes a dig, someone else needs to set" $null $null

# Row 14 - new row with wrapped, multi-line details and a taller row to match
Set-Row 14 "Use the distance away from the ball to effect the ability to make a dig or pass" "Garrett" "The farther a player is away from the ball, the less likely they are to make a dig. Going to try to have it feel lik 50% is player skill and 50% is proximity to the ball`nequal chance to dig a hard hit right at the target as digging an average hit 3-4 squares away`nNeed to consider that the y axis might logically have a greater impact than the x axis to accound for forearm and overhandhand digs`nmaybe a system where one tile in the y axis adds 1 difficulty and one tile in the x axis adds 1/2 difficulty" $true
$ws.Rows.Item(14).RowHeight = 57.6

# Row 15 - new row
Set-Row 15 "BUG -AI not always going to the correct locations off digs" "Garrett" "Fix it "

# ---- Tobi section (rows 19-23) ----
Set-Row 19 "BUG - Some squares full" "Tobi" "Something happens during the rally and the player can't move their pieces to certain tiles"
Set-Row 20 "Limit movement during player blocker/defender reaction phase" "Tobi" "during the player blocker reaction phase all should be able to move, and their movement limited to one square from their starting location"

Set-Row 21 "Functionality for determining how many blocks in the area" "Tobi" "When a pawn is attacking, determine how many opponents blocks are in a valid area to add to the block quality value of the simulation`nwould also be fantastic to take into consideration where the ball is headed when calculating which blocks are close enough to be a factor" $true
$ws.Rows.Item(21).RowHeight = 28.8

# Row 22 used to hold the "closest player" feature row - content moved up into
# the Garrett section above, so this row is now blank (formatting untouched).
$ws.Range("B22:D22").ClearContents()
$ws.Rows.Item(22).AutoFit() | Out-Null

Set-Row 23 "Turn off setters movement only during certain phases" "Tobi" "There are certain rally phases where I'd like to set the setters position based on the ball position and not let the player move the setter from there at all"

# ---- Sheet-level tweaks ----
# NOTE: the engine's ColumnWidth setter quantizes to a 1/6-character grid
# (mirrors Excel's own pixel-grid rounding); 63.17 is the input that round-
# trips to the target stored width of 64.
$ws.Columns.Item(2).ColumnWidth = 63.17
$ws.Range("B15").Select() | Out-Null
